$d = $word.ActiveDocument

$pairs = @(
    @("676÷6=112, 4", "326÷2=163, 0"),
    @("256÷3=85, 1",  "469÷9=52, 1"),
    @("602÷6=100, 2", "777÷6=129, 3"),
    @("200÷3=66, 2",  "439÷6=73, 1"),
    @("462÷8=57, 6",  "692÷5=138, 2"),
    @("253÷9=28, 1",  "146÷8=18, 2"),
    @("452÷9=50, 2",  "492÷6=82, 0"),
    @("630÷4=157, 2", "835÷4=208, 3"),
    @("520÷5=104, 0", "109÷2=54, 1"),
    @("824÷5=164, 4", "692÷2=346, 0"),
    @("139÷7=19, 6",  "188÷9=20, 8"),
    @("430÷8=53, 6",  "169÷9=18, 7"),
    @("333÷7=47, 4",  "359÷5=71, 4"),
    @("286÷8=35, 6",  "329÷6=54, 5"),
    @("949÷9=105, 4", "990÷3=330, 0"),
    @("441÷4=110, 1", "220÷4=55, 0"),
    @("859÷4=214, 3", "935÷2=467, 1"),
    @("160÷7=22, 6",  "819÷8=102, 3"),
    @("648÷9=72, 0",  "479÷4=119, 3"),
    @("244÷4=61, 0",  "577÷5=115, 2"),
    @("500÷3=166, 2", "434÷8=54, 2"),
    @("484÷4=121, 0", "322÷6=53, 4"),
    @("656÷7=93, 5",  "869÷3=289, 2"),
    @("589÷2=294, 1", "269÷7=38, 3"),
    @("479÷5=95, 4",  "141÷4=35, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
